# Apply edits to "Monitoramento diário" worksheet (first sheet of the workbook)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monitoramento diário")

# Update the daily figures for row 15 (B15 = day 22)
$ws.Range("C15").Value = 5
$ws.Range("E15").Value = 43
$ws.Range("G15").Value = 43
$ws.Range("I15").Value = 43

# Move the active selection to J14 (matches saved sheetView selection)
$ws.Activate()
$ws.Range("J14").Select()
